# Adds a totals row (row 15) with SUM formulas to both sheets, then
# leaves the "Single" sheet active/selected (matching the author's last
# on-screen state): B15:M15 selected on "Group", J17 selected on "Single".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Group")
$ws2 = $wb.Worksheets.Item("Single")

foreach ($ws in @($ws1, $ws2)) {
    $ws.Range("B15").Formula = "=SUM(B2:B14)"
    $ws.Range("C15:M15").Formula = "=SUM(C2:C14)"
}

# Page setup tweak recorded on the "Group" sheet only.
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# Final on-screen state: "Single" is the active/visible tab, with its own
# selection; "Group" keeps a remembered selection over the new totals row.
$ws1.Activate() | Out-Null
$ws1.Range("B15:M15").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("J17").Select() | Out-Null
